# Refresh NATMI LR-pair (Ptdss1 -> Jmjd6) TPM-derived metrics.
# The underlying TPM recomputation changes ligand/receptor average & total
# expression values (columns G, H, M, N), their derived-specificity scores
# (I, J, O, P), and the resulting edge-weight / edge-specificity scores
# (Q, R, S, T) for every Sending-cluster x Target-cluster row (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.176422
$ws.Range("H2").Value = 18.529266
$ws.Range("I2").Value = 0.3058063741187975
$ws.Range("J2").Value = 0.3058063741187975
$ws.Range("M2").Value = 7.446634
$ws.Range("N2").Value = 22.339902
$ws.Range("O2").Value = 0.3025192762093004
$ws.Range("P2").Value = 0.3025192762093004
$ws.Range("Q2").Value = 45.993554063548
$ws.Range("R2").Value = 413.941986571932
$ws.Range("S2").Value = 0.09251232295860916
$ws.Range("T2").Value = 0.09251232295860916
$ws.Range("G3").Value = 6.176422
$ws.Range("H3").Value = 18.529266
$ws.Range("I3").Value = 0.3058063741187975
$ws.Range("J3").Value = 0.3058063741187975
$ws.Range("O3").Value = 0.4584869690672005
$ws.Range("P3").Value = 0.4584869690672005
$ws.Range("Q3").Value = 69.70612075851665
$ws.Range("R3").Value = 627.3550868266499
$ws.Range("S3").Value = 0.1402082375911579
$ws.Range("T3").Value = 0.1402082375911579
$ws.Range("G4").Value = 6.176422
$ws.Range("H4").Value = 18.529266
$ws.Range("I4").Value = 0.3058063741187975
$ws.Range("J4").Value = 0.3058063741187975
$ws.Range("M4").Value = 5.882927666666666
$ws.Range("O4").Value = 0.238993754723499
$ws.Range("P4").Value = 0.238993754723499
$ws.Range("Q4").Value = 36.33544386480866
$ws.Range("R4").Value = 327.0189947832779
$ws.Range("S4").Value = 0.07308581356903047
$ws.Range("T4").Value = 0.07308581356903046
$ws.Range("I5").Value = 0.4631823009753332
$ws.Range("J5").Value = 0.4631823009753332
$ws.Range("M5").Value = 7.446634
$ws.Range("N5").Value = 22.339902
$ws.Range("O5").Value = 0.3025192762093004
$ws.Range("P5").Value = 0.3025192762093004
$ws.Range("Q5").Value = 69.66303518876867
$ws.Range("R5").Value = 626.967316698918
$ws.Range("S5").Value = 0.1401215744440162
$ws.Range("T5").Value = 0.1401215744440161
$ws.Range("I6").Value = 0.4631823009753332
$ws.Range("J6").Value = 0.4631823009753332
$ws.Range("O6").Value = 0.4584869690672005
$ws.Range("P6").Value = 0.4584869690672005
$ws.Range("S6").Value = 0.2123630492997524
$ws.Range("T6").Value = 0.2123630492997524
$ws.Range("I7").Value = 0.4631823009753332
$ws.Range("J7").Value = 0.4631823009753332
$ws.Range("M7").Value = 5.882927666666666
$ws.Range("O7").Value = 0.238993754723499
$ws.Range("P7").Value = 0.238993754723499
$ws.Range("S7").Value = 0.1106976772315647
$ws.Range("T7").Value = 0.1106976772315647
$ws.Range("I8").Value = 0.2310113249058692
$ws.Range("J8").Value = 0.2310113249058692
$ws.Range("M8").Value = 7.446634
$ws.Range("N8").Value = 22.339902
$ws.Range("O8").Value = 0.3025192762093004
$ws.Range("P8").Value = 0.3025192762093004
$ws.Range("Q8").Value = 34.744311304716
$ws.Range("R8").Value = 312.6988017424441
$ws.Range("S8").Value = 0.06988537880667509
$ws.Range("T8").Value = 0.06988537880667509
$ws.Range("I9").Value = 0.2310113249058692
$ws.Range("J9").Value = 0.2310113249058692
$ws.Range("O9").Value = 0.4584869690672005
$ws.Range("P9").Value = 0.4584869690672005
$ws.Range("S9").Value = 0.1059156821762903
$ws.Range("T9").Value = 0.1059156821762903
$ws.Range("I10").Value = 0.2310113249058692
$ws.Range("J10").Value = 0.2310113249058692
$ws.Range("M10").Value = 5.882927666666666
$ws.Range("O10").Value = 0.238993754723499
$ws.Range("P10").Value = 0.238993754723499
$ws.Range("S10").Value = 0.05521026392290385
$ws.Range("T10").Value = 0.05521026392290384
